# Weekly refresh of Fruta/Vega Modelo de Temuco - Maracuyá rows: each destination
# row (2-34) is repopulated with the Fecha/Volumen/Precio/Unidad/Origen/Kg data
# that used to live on a different row (source records were re-sorted by date).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 12)
$ws.Range("D2").Value = 44379
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 30000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 30000
$ws.Range("Q2").Value = "`$/caja 18 kilos"
$ws.Range("R2").Value = "Región de Arica y Parinacota"
$ws.Range("S2").Value = 1667
$ws.Range("T2").Value = 18

# Row 3 (was row 19)
$ws.Range("D3").Value = 44671
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 32000
$ws.Range("O3").Value = 32000
$ws.Range("P3").Value = 32000
$ws.Range("Q3").Value = "`$/caja 18 kilos"
$ws.Range("R3").Value = "Región de Arica y Parinacota"
$ws.Range("S3").Value = 1778
$ws.Range("T3").Value = 18

# Row 4 (was row 29)
$ws.Range("D4").Value = 44634
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 45000
$ws.Range("O4").Value = 45000
$ws.Range("P4").Value = 45000
$ws.Range("Q4").Value = "`$/caja 18 kilos"
$ws.Range("R4").Value = "Región de Arica y Parinacota"
$ws.Range("S4").Value = 2500
$ws.Range("T4").Value = 18

# Row 5 (was row 14)
$ws.Range("D5").Value = 44662
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("Q5").Value = "`$/caja 18 kilos"
$ws.Range("R5").Value = "Región de Arica y Parinacota"
$ws.Range("S5").Value = 1667
$ws.Range("T5").Value = 18

# Row 6 (was row 26)
$ws.Range("D6").Value = 44645
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("Q6").Value = "`$/caja 18 kilos"
$ws.Range("R6").Value = "Región de Arica y Parinacota"
$ws.Range("S6").Value = 1667
$ws.Range("T6").Value = 18

# Row 7 (was row 31)
$ws.Range("D7").Value = 44432
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 35000
$ws.Range("O7").Value = 35000
$ws.Range("P7").Value = 35000
$ws.Range("Q7").Value = "`$/caja 18 kilos"
$ws.Range("R7").Value = "Perú"
$ws.Range("S7").Value = 1944
$ws.Range("T7").Value = 18

# Row 8 (was row 16)
$ws.Range("D8").Value = 44433
$ws.Range("M8").Value = 15
$ws.Range("N8").Value = 35000
$ws.Range("O8").Value = 35000
$ws.Range("P8").Value = 35000
$ws.Range("Q8").Value = "`$/caja 18 kilos"
$ws.Range("R8").Value = "Región de Arica y Parinacota"
$ws.Range("S8").Value = 1944
$ws.Range("T8").Value = 18

# Row 9 (was row 32)
$ws.Range("D9").Value = 44438
$ws.Range("M9").Value = 25
$ws.Range("N9").Value = 35000
$ws.Range("O9").Value = 35000
$ws.Range("P9").Value = 35000
$ws.Range("Q9").Value = "`$/caja 18 kilos"
$ws.Range("R9").Value = "Región de Arica y Parinacota"
$ws.Range("S9").Value = 1944
$ws.Range("T9").Value = 18

# Row 10 (was row 22)
$ws.Range("D10").Value = 44669
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 32000
$ws.Range("O10").Value = 32000
$ws.Range("P10").Value = 32000
$ws.Range("Q10").Value = "`$/caja 18 kilos"
$ws.Range("R10").Value = "Región de Arica y Parinacota"
$ws.Range("S10").Value = 1778
$ws.Range("T10").Value = 18

# Row 11 (was row 23)
$ws.Range("D11").Value = 44363
$ws.Range("M11").Value = 144
$ws.Range("N11").Value = 1700
$ws.Range("O11").Value = 1700
$ws.Range("P11").Value = 1700
$ws.Range("Q11").Value = "`$/kilo"
$ws.Range("R11").Value = "Región de Arica y Parinacota"
$ws.Range("S11").Value = 1700
$ws.Range("T11").Value = 1

# Row 12 (was row 9)
$ws.Range("D12").Value = 44676
$ws.Range("M12").Value = 55
$ws.Range("N12").Value = 28000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 28909
$ws.Range("Q12").Value = "`$/caja 18 kilos"
$ws.Range("R12").Value = "Región de Arica y Parinacota"
$ws.Range("S12").Value = 1606
$ws.Range("T12").Value = 18

# Row 13 (was row 3)
$ws.Range("D13").Value = 44629
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 35000
$ws.Range("O13").Value = 35000
$ws.Range("P13").Value = 35000
$ws.Range("Q13").Value = "`$/caja 18 kilos"
$ws.Range("R13").Value = "Región de Arica y Parinacota"
$ws.Range("S13").Value = 1944
$ws.Range("T13").Value = 18

# Row 14 (was row 15)
$ws.Range("D14").Value = 44392
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = 35000
$ws.Range("O14").Value = 35000
$ws.Range("P14").Value = 35000
$ws.Range("Q14").Value = "`$/caja 18 kilos"
$ws.Range("R14").Value = "Región de Arica y Parinacota"
$ws.Range("S14").Value = 1944
$ws.Range("T14").Value = 18

# Row 15 (was row 21)
$ws.Range("D15").Value = 44658
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 28000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 28000
$ws.Range("Q15").Value = "`$/caja 18 kilos"
$ws.Range("R15").Value = "Región de Arica y Parinacota"
$ws.Range("S15").Value = 1556
$ws.Range("T15").Value = 18

# Row 16 (was row 30)
$ws.Range("D16").Value = 44385
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 1200
$ws.Range("O16").Value = 1200
$ws.Range("P16").Value = 1200
$ws.Range("Q16").Value = "`$/kilo"
$ws.Range("R16").Value = "Perú"
$ws.Range("S16").Value = 1200
$ws.Range("T16").Value = 1

# Row 17 (was row 8)
$ws.Range("D17").Value = 44264
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 40000
$ws.Range("O17").Value = 40000
$ws.Range("P17").Value = 40000
$ws.Range("Q17").Value = "`$/caja 18 kilos"
$ws.Range("R17").Value = "Región de Arica y Parinacota"
$ws.Range("S17").Value = 2222
$ws.Range("T17").Value = 18

# Row 18 (was row 6)
$ws.Range("D18").Value = 44294
$ws.Range("M18").Value = 15
$ws.Range("N18").Value = 35000
$ws.Range("O18").Value = 35000
$ws.Range("P18").Value = 35000
$ws.Range("Q18").Value = "`$/caja 18 kilos"
$ws.Range("R18").Value = "Región de Arica y Parinacota"
$ws.Range("S18").Value = 1944
$ws.Range("T18").Value = 18

# Row 19 (was row 24)
$ws.Range("D19").Value = 44424
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = 35000
$ws.Range("O19").Value = 35000
$ws.Range("P19").Value = 35000
$ws.Range("Q19").Value = "`$/caja 18 kilos"
$ws.Range("R19").Value = "Región de Arica y Parinacota"
$ws.Range("S19").Value = 1944
$ws.Range("T19").Value = 18

# Row 20 (was row 10)
$ws.Range("D20").Value = 44279
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = 35000
$ws.Range("O20").Value = 36000
$ws.Range("P20").Value = 35667
$ws.Range("Q20").Value = "`$/caja 18 kilos"
$ws.Range("R20").Value = "Región de Arica y Parinacota"
$ws.Range("S20").Value = 1982
$ws.Range("T20").Value = 18

# Row 21 (was row 28)
$ws.Range("D21").Value = 44434
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 35000
$ws.Range("O21").Value = 35000
$ws.Range("P21").Value = 35000
$ws.Range("Q21").Value = "`$/caja 18 kilos"
$ws.Range("R21").Value = "Región de Arica y Parinacota"
$ws.Range("S21").Value = 1944
$ws.Range("T21").Value = 18

# Row 22 (was row 17)
$ws.Range("D22").Value = 44377
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 40000
$ws.Range("O22").Value = 40000
$ws.Range("P22").Value = 40000
$ws.Range("Q22").Value = "`$/caja 18 kilos"
$ws.Range("R22").Value = "Región de Arica y Parinacota"
$ws.Range("S22").Value = 2222
$ws.Range("T22").Value = 18

# Row 23 (was row 20)
$ws.Range("D23").Value = 44664
$ws.Range("M23").Value = 15
$ws.Range("N23").Value = 30000
$ws.Range("O23").Value = 30000
$ws.Range("P23").Value = 30000
$ws.Range("Q23").Value = "`$/caja 18 kilos"
$ws.Range("R23").Value = "Perú"
$ws.Range("S23").Value = 1667
$ws.Range("T23").Value = 18

# Row 24 (was row 5)
$ws.Range("D24").Value = 44364
$ws.Range("M24").Value = 90
$ws.Range("N24").Value = 1700
$ws.Range("O24").Value = 1700
$ws.Range("P24").Value = 1700
$ws.Range("Q24").Value = "`$/kilo"
$ws.Range("R24").Value = "Región de Arica y Parinacota"
$ws.Range("S24").Value = 1700
$ws.Range("T24").Value = 1

# Row 25 (was row 27)
$ws.Range("D25").Value = 44405
$ws.Range("M25").Value = 10
$ws.Range("N25").Value = 35000
$ws.Range("O25").Value = 35000
$ws.Range("P25").Value = 35000
$ws.Range("Q25").Value = "`$/caja 18 kilos"
$ws.Range("R25").Value = "Región de Arica y Parinacota"
$ws.Range("S25").Value = 1944
$ws.Range("T25").Value = 18

# Row 26 (was row 33)
$ws.Range("D26").Value = 44435
$ws.Range("M26").Value = 10
$ws.Range("N26").Value = 35000
$ws.Range("O26").Value = 35000
$ws.Range("P26").Value = 35000
$ws.Range("Q26").Value = "`$/caja 18 kilos"
$ws.Range("R26").Value = "Perú"
$ws.Range("S26").Value = 1944
$ws.Range("T26").Value = 18

# Row 27 (was row 34)
$ws.Range("D27").Value = 44435
$ws.Range("M27").Value = 105
$ws.Range("N27").Value = 35000
$ws.Range("O27").Value = 35000
$ws.Range("P27").Value = 35000
$ws.Range("Q27").Value = "`$/caja 18 kilos"
$ws.Range("R27").Value = "Región de Arica y Parinacota"
$ws.Range("S27").Value = 1944
$ws.Range("T27").Value = 18

# Row 28 (was row 25)
$ws.Range("D28").Value = 44442
$ws.Range("M28").Value = 15
$ws.Range("N28").Value = 35000
$ws.Range("O28").Value = 35000
$ws.Range("P28").Value = 35000
$ws.Range("Q28").Value = "`$/caja 18 kilos"
$ws.Range("R28").Value = "Perú"
$ws.Range("S28").Value = 1944
$ws.Range("T28").Value = 18

# Row 29 (was row 7)
$ws.Range("D29").Value = 44418
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 35000
$ws.Range("O29").Value = 35000
$ws.Range("P29").Value = 35000
$ws.Range("Q29").Value = "`$/caja 18 kilos"
$ws.Range("R29").Value = "Región de Arica y Parinacota"
$ws.Range("S29").Value = 1944
$ws.Range("T29").Value = 18

# Row 30 (was row 13)
$ws.Range("D30").Value = 44449
$ws.Range("M30").Value = 20
$ws.Range("N30").Value = 38000
$ws.Range("O30").Value = 38000
$ws.Range("P30").Value = 38000
$ws.Range("Q30").Value = "`$/caja 18 kilos"
$ws.Range("R30").Value = "Región de Arica y Parinacota"
$ws.Range("S30").Value = 2111
$ws.Range("T30").Value = 18

# Row 31 (was row 11)
$ws.Range("D31").Value = 44448
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 38000
$ws.Range("O31").Value = 38000
$ws.Range("P31").Value = 38000
$ws.Range("Q31").Value = "`$/caja 18 kilos"
$ws.Range("R31").Value = "Región de Arica y Parinacota"
$ws.Range("S31").Value = 2111
$ws.Range("T31").Value = 18

# Row 32 (was row 4)
$ws.Range("D32").Value = 44431
$ws.Range("M32").Value = 30
$ws.Range("N32").Value = 35000
$ws.Range("O32").Value = 35000
$ws.Range("P32").Value = 35000
$ws.Range("Q32").Value = "`$/caja 18 kilos"
$ws.Range("R32").Value = "Región de Arica y Parinacota"
$ws.Range("S32").Value = 1944
$ws.Range("T32").Value = 18

# Row 33 (was row 2)
$ws.Range("D33").Value = 44357
$ws.Range("M33").Value = 10
$ws.Range("N33").Value = 38000
$ws.Range("O33").Value = 38000
$ws.Range("P33").Value = 38000
$ws.Range("Q33").Value = "`$/caja 18 kilos"
$ws.Range("R33").Value = "Perú"
$ws.Range("S33").Value = 2111
$ws.Range("T33").Value = 18

# Row 34 (was row 18)
$ws.Range("D34").Value = 44369
$ws.Range("M34").Value = 5
$ws.Range("N34").Value = 35000
$ws.Range("O34").Value = 35000
$ws.Range("P34").Value = 35000
$ws.Range("Q34").Value = "`$/caja 18 kilos"
$ws.Range("R34").Value = "Perú"
$ws.Range("S34").Value = 1944
$ws.Range("T34").Value = 18
